$d = $word.ActiveDocument

# Unicode helpers
$ENDASH = [char]0x2013
$RSQUOTE = [char]0x2019

# ---------------------------------------------------------------------------
# 1) Title: "Board Game Work Log" -> add bold " -- Stephanie" after it.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1   # exclude the paragraph mark
$r1.Collapse(0)         # collapse to end (wdCollapseEnd)
$r1.InsertAfter(" -- Stephanie")
$r1.Font.Bold = -1

# ---------------------------------------------------------------------------
# 2) September 28th paragraph: split sentence + add GitHub parenthetical.
# ---------------------------------------------------------------------------
$old2 = " (Stephanie) $ENDASH I had done some little edits and putzing around that I did not record here in previous days. Today I changed GUI"
$new2 = " (Stephanie) $ENDASH I had done some little edits and putzing around that I did not record here in previous days (should all be reflected in Github${RSQUOTE}s commits though). Today I changed GUI"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) September 30th paragraph: rewrite to "(All)" with extra text.
# ---------------------------------------------------------------------------
$old3 = " (Stephanie) $ENDASH Group meeting where we worked on the code together."
$new3 = " (All) $ENDASH There have been other meetings that I didn${RSQUOTE}t jot down here. Group meeting where we worked on the code together today. Lots of testing."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Insert a new "October 2nd (All) ..." paragraph + a blank eastAsia-hint
#    paragraph right before the "October 3rd" paragraph.
# ---------------------------------------------------------------------------
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$findOct3 = $d.Content
$findOct3.Find.Execute("October 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$oct3Start = $findOct3.Start
$insertPoint = $d.Range($oct3Start, $oct3Start)

$oct2Xml = "<w:p $ns><w:r><w:t>October 2</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=`"superscript`"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:t xml:space=`"preserve`"> (All) $ENDASH Group meeting where we did lots of testing. We also compiled all our times together and </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>filled</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> out the Team Tracking sheet. </w:t></w:r></w:p><w:p $ns><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr></w:p>"
$insertPoint.InsertXML($oct2Xml)

# ---------------------------------------------------------------------------
# 5) October 3rd paragraph: split lead-in, drop the "Once again..." phrase.
# ---------------------------------------------------------------------------
$old5 = " (Stephanie) $ENDASH Once again, there have been meetings that I didn${RSQUOTE}t jot down here. Today I mostly cleaned up the comment section of GUI and made it so the list of numbers in your hand wrap when they reach the end of the screen."
$new5 = " (Stephanie) $ENDASH Today I mostly cleaned up the comment section of GUI and made it so the list of numbers in your hand wrap when they reach the end of the screen."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Append a blank paragraph + new "October 4th (Stephanie) ..." paragraph
#    at the very end of the document.
# ---------------------------------------------------------------------------
$endR = $d.Range($d.Content.End, $d.Content.End)
$oct4Xml = "<w:p $ns/><w:p $ns><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:t>October 4</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=`"superscript`"/></w:rPr><w:t>th</w:t></w:r><w:r><w:t xml:space=`"preserve`"> (Stephanie) $ENDASH Small cleaning up of code and comments again for the GUI.</w:t></w:r></w:p>"
$endR.InsertXML($oct4Xml)

Write-Output "Done"
